$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tests")

# --- Row 11 : Alimentation / résistance test now has a result + tightened wording ---
$ws.Range("C11").Value = "Tester la résistance entre les pattes des composants énergivores avec un multimètre en mode résistance sur les pattes 5v et la patte des pattes des Vin des composantes pour éviter les court-circuit."
$ws.Range("D11").Value = "Réussi"

# --- Row 14 : Connecteur USB-C now Réussi, comment reworded ("boucle" instead of "petit loop") ---
$ws.Range("D14").Value = "Réussi"
$ws.Range("E14").Value = "la longueur prévue pour le fils USB-C est un peu court. Nous allons devoir faire une boucle"

# --- Row 17 : Position des connecteurs -> Réussi ---
$ws.Range("D17").Value = "Réussi"

# --- Row 18 : Alimentation des esp32 -> Réussi ---
$ws.Range("D18").Value = "Réussi"

# --- Row 19 : Communication i2c -> Échec, new comment ---
$ws.Range("D19").Value = "Échec"
$ws.Range("E19").Value = "Il faut reset les esp32 pour activer le i2c des esp32"

# --- Row 18 extra: E18 picks up a (still-empty) centered/wrapped white-fill style ---
$ws.Range("E18").Interior.ThemeColor = 2
$ws.Range("E18").Interior.TintAndShade = 0

# --- Row 20 : Mini écran -> Réussi, taller row to match new wrapped content nearby ---
$ws.Range("D20").Value = "Réussi"
$ws.Rows(20).RowHeight = 27.6

# --- Row 21 : Alimentation intégrale -> Échec, new comment, taller row ---
$ws.Range("D21").Value = "Échec"
$ws.Range("E21").Value = "L'écran demande trop de courant si on l'alimente du PI"
$ws.Rows(21).RowHeight = 52.2

# --- Bottom "modifications / comments" table: row 34 comment expanded ---
$ws.Range("C34").Value = "Le pi n'a pas assez de courant pour alimenter tous et l'écran avec les ports USB-C du pi"

# --- Restore the view: scrolled up a bit and the window maximised on save ---
$ws.Application.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.WindowState = -4137
